$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "maa://24702 (94.48), maa://25390 (96.14), maa://36681 (87.34)"
$ws.Range("L2").Value = "*maa://24633 (56.52), *maa://30515 (69.9), *maa://34787 (72.97), maa://39402 (91.23), ***maa://20792 (11.93), ***maa://29083 (27.78)"
$ws.Range("D3").Value = "maa://36987 (96.08), maa://40192 (100.0), maa://39849 (88.89)"
$ws.Range("T3").Value = "maa://24617 (89.74), **maa://20790 (43.48), ***maa://37170 (16.92), maa://45854 (100.0)"
$ws.Range("AF3").Value = "*maa://21289 (75.0)"
$ws.Range("D4").Value = "maa://24632 (94.01), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("T4").Value = "maa://32509 (96.52), maa://27295 (85.71), maa://22754 (90.41), *maa://21746 (55.81), *maa://31008 (78.57)"
$ws.Range("X4").Value = "**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (88.68), ***maa://36683 (28.26)"
$ws.Range("D5").Value = "maa://21245 (84.68), maa://22744 (84.0)"
$ws.Range("D7").Value = "maa://21955 (94.87)"
$ws.Range("L7").Value = "maa://28624 (92.52), maa://24957 (97.73)"
$ws.Range("X7").Value = "maa://22399 (95.33), *maa://22758 (74.29)"
$ws.Range("A8").Value = "更新日期：2025.02.12 13:18:17"
$ws.Range("D8").Value = "*maa://21476 (73.08), *maa://39431 (57.14), *maa://37551 (57.14)"
$ws.Range("X8").Value = "maa://21411 (95.88)"
$ws.Range("L9").Value = "maa://22762 (92.22), *maa://39552 (75.0)"
$ws.Range("X9").Value = "maa://26223 (97.83)"
$ws.Range("AB9").Value = "maa://28711 (87.18), ***maa://22740 (5.66), **maa://39938 (46.67), **maa://27377 (42.86), ***maa://25174 (19.05), maa://40166 (96.0), *maa://45044 (66.67)"
$ws.Range("D10").Value = "***maa://25695 (18.82), ***maa://34206 (20.0), ***maa://39951 (15.69), ***maa://39243 (28.57), *maa://45271 (57.69)"
$ws.Range("AF10").Value = "*maa://25021 (53.85), *maa://22733 (60.0), **maa://22761 (50.0)"
$ws.Range("T11").Value = "maa://22747 (93.08), maa://22501 (97.67), *maa://45521 (76.92)"
$ws.Range("X11").Value = "maa://36713 (97.74)"
$ws.Range("H12").Value = "maa://21867 (89.88), ***maa://45826 (25.0)"
$ws.Range("X12").Value = "maa://22753 (90.91), *maa://21485 (76.43), maa://37962 (90.24)"
$ws.Range("AB12").Value = "maa://23669 (95.49), maa://36677 (93.1), maa://39872 (91.3)"
$ws.Range("AF12").Value = "*maa://28932 (77.93), *maa://20106 (63.96), *maa://22769 (64.29)"
$ws.Range("D13").Value = "maa://24999 (92.03), maa://36673 (93.15), maa://25001 (85.71)"
$ws.Range("H13").Value = "*maa://21248 (73.19), **maa://22728 (47.73)"
$ws.Range("AF13").Value = "**maa://22737 (33.33), maa://39883 (91.18), *maa://39885 (53.33)"
$ws.Range("T14").Value = "maa://22521 (94.29), maa://42751 (100.0)"
$ws.Range("X14").Value = "maa://37468 (90.91)"
$ws.Range("D15").Value = "*maa://22743 (77.62), maa://22734 (84.17), *maa://30808 (64.18), **maa://36048 (44.07), maa://45058 (91.67)"
$ws.Range("X15").Value = "maa://38786 (85.71)"
$ws.Range("AF15").Value = "maa://21364 (81.23), *maa://36666 (79.25), *maa://22766 (68.64)"
$ws.Range("D16").Value = "maa://21441 (96.4), maa://36679 (94.23), maa://37650 (97.14)"
$ws.Range("T16").Value = "maa://22729 (94.94), *maa://28648 (69.12), maa://36674 (81.25)"
$ws.Range("H17").Value = "maa://22430 (88.66), maa://39599 (85.71)"
$ws.Range("D18").Value = "maa://24570 (97.3)"
$ws.Range("L18").Value = "maa://22466 (89.94), *maa://22732 (51.14)"
$ws.Range("D20").Value = "maa://21432 (90.42), maa://25198 (93.58), *maa://20795 (51.16), maa://36680 (93.75)"
$ws.Range("L20").Value = "maa://41331 (85.52)"
$ws.Range("D21").Value = "maa://21261 (97.5)"
$ws.Range("H21").Value = "maa://24372 (96.94)"
$ws.Range("X21").Value = "maa://20110 (86.76), maa://34946 (93.02)"
$ws.Range("AB21").Value = "maa://21443 (80.97), ***maa://23820 (30.0)"
$ws.Range("AF21").Value = "maa://22524 (94.5), *maa://22432 (77.78)"
$ws.Range("X22").Value = "maa://21282 (98.59), *maa://37649 (65.52)"
$ws.Range("L23").Value = "maa://39756 (95.51), maa://39875 (94.37)"
$ws.Range("X23").Value = "*maa://28503 (68.35)"
$ws.Range("D24").Value = "*maa://24368 (78.16), **maa://46650 (50.0)"
$ws.Range("X24").Value = "maa://29988 (84.92), maa://23504 (93.08), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.63), ***maa://22815 (23.08)"
$ws.Range("AF24").Value = "maa://22523 (85.93), maa://36672 (80.36), maa://29910 (92.98), **maa://21440 (35.71), *maa://45831 (75.0)"
$ws.Range("D25").Value = "maa://29753 (95.09)"
$ws.Range("H25").Value = "*maa://29063 (74.21), *maa://25311 (73.53), ***maa://22725 (4.84), *maa://45047 (62.5)"
$ws.Range("AB26").Value = "maa://42235 (94.74)"
$ws.Range("H27").Value = "**maa://21283 (47.37), *maa://39601 (78.95), maa://34494 (97.14), **maa://36665 (50.0)"
$ws.Range("T27").Value = "*maa://30624 (78.33)"
$ws.Range("D28").Value = "maa://24465 (91.1), maa://25725 (83.72)"
$ws.Range("X28").Value = "maa://39929 (90.53), maa://41749 (90.36), ***maa://39723 (13.89)"
$ws.Range("L29").Value = "maa://28432 (93.39), *maa://28440 (79.44), maa://31400 (100.0), *maa://28650 (71.43)"
$ws.Range("AF29").Value = "*maa://24080 (68.77), maa://42865 (81.03), ***maa://34960 (8.33)"
$ws.Range("L30").Value = "maa://30442 (95.24)"
$ws.Range("AB30").Value = "maa://42979 (96.77), maa://45822 (100.0), *maa://45045 (80.0)"
$ws.Range("L31").Value = "maa://35926 (93.36), maa://36258 (84.96), *maa://43904 (72.73)"
$ws.Range("H32").Value = "maa://21895 (97.5), maa://36667 (98.72), **maa://20793 (38.78), maa://22760 (100.0)"
$ws.Range("T32").Value = "maa://42859 (96.52), maa://41108 (88.0), maa://41238 (97.09), maa://45523 (100.0)"
$ws.Range("L37").Value = "maa://45718 (99.11), maa://45789 (100.0)"
$ws.Range("AF38").Value = "maa://36697 (86.06)"
$ws.Range("H39").Value = "maa://36670 (88.89), maa://25199 (84.82), maa://30434 (91.14), ***maa://25036 (16.0), maa://45059 (81.25), *maa://44165 (66.67)"
$ws.Range("P39").Value = "maa://24709 (91.33)"
$ws.Range("T39").Value = "maa://45788 (82.02), maa://45790 (81.82)"
$ws.Range("P40").Value = "maa://23278 (95.53), maa://21386 (95.77), maa://36664 (89.29), maa://45550 (100.0)"
$ws.Range("H45").Value = "maa://21229 (84.66), maa://30807 (95.65), *maa://22767 (55.0), ***maa://20796 (13.79), maa://42459 (84.21)"
$ws.Range("T45").Value = "**maa://39364 (36.67)"
$ws.Range("H46").Value = "maa://35931 (92.63), maa://43901 (91.67)"
$ws.Range("H47").Value = "maa://27410 (96.43), maa://29661 (97.3), maa://28038 (84.62)"
$ws.Range("H53").Value = "maa://32534 (93.9), **maa://32434 (33.33)"
$ws.Range("H55").Value = "maa://32532 (92.23)"
$ws.Range("H59").Value = "maa://31270 (95.28), maa://27746 (82.3)"
$ws.Range("H60").Value = "*maa://40438 (67.8)"
$ws.Range("H62").Value = "maa://42981 (94.87), maa://43903 (100.0)"
